# Insert a new weekly data row above row 33 (shifts existing rows 33:83 down to 34:84)
# and populate the new row 33 with the latest week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 33:83 down to 34:84 by inserting a new blank row at 33.
$ws.Rows.Item(33).Insert()

# The insert operation copies formatting from the row above (row 32), which already
# matches the rest of the table, so we only need to set the cell values.
$ws.Cells.Item(33, 1).Value = 7
$ws.Cells.Item(33, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(33, 3).Value = "Ñuble"
$ws.Cells.Item(33, 4).Value = 44662
$ws.Cells.Item(33, 5).Value = 16
$ws.Cells.Item(33, 6).Value = 100112031
$ws.Cells.Item(33, 7).Value = "Poroto verde"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 120
$ws.Cells.Item(33, 11).Value = 24000
$ws.Cells.Item(33, 12).Value = 25000
$ws.Cells.Item(33, 13).Value = 24500
$ws.Cells.Item(33, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(33, 16).Value = 980
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
